# Loan RBI, Variable Instalments
# Adds a new (currently blank) column to the "Repayment schedule" sheet,
# right before the existing "Late" column, shifting "Late" / "Outstanding" /
# heading one column to the right, and switches the active sheet/selection
# to "Repayment schedule".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14) - shifts N,O,P -> O,P,Q
$ws.Columns.Item(14).EntireColumn.Insert()

# The newly inserted column keeps the default sheet width; give it the
# same width the author set when adding the new data column.
$ws.Columns.Item(14).ColumnWidth = 10.1666666666667

# Make "Repayment schedule" the active sheet/tab with the new selection
$ws.Activate()
$ws.Range("S9").Select()
